# Port-level landings table: move the "EUREKA AREA TOTALS" label from column B
# to column A on row 2, and put a new "Totals" label into column B (row 2)
# with default/normal formatting. Also refresh the column widths that
# auto-resized as a result, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the label currently sitting in B2 ("EUREKA AREA TOTALS ")
$oldB2Value = $ws.Range("B2").Value()

# Carry B2's formatting (font, etc.) over to A2, then carry the text itself
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A2").Value = $oldB2Value

# B2 becomes a plain "Totals" cell with default (Normal) styling
$ws.Range("B2").Style = "Normal"
$ws.Range("B2").Value = "Totals"

# Columns A and B resize (bestFit) to match their new longest contents
$ws.Columns(1).ColumnWidth = 19.498697916666668
$ws.Columns(2).ColumnWidth = 17.330729166666668

# Selection moves to A2, and the view scrolls back so A1 is the top-left cell
$ws.Range("A2").Select()
